$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Phone-number table cell (row 4, col 4): "0169 777 1714" -> "0166 907 0077"
#    split across two runs with a _GoBack bookmark sitting between them.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)
$cell = $t.Cell(4, 4)
$cellRng = $cell.Range
$textRng = $d.Range($cellRng.Start, $cellRng.End - 1)
$textRng.Text = "0166 907 0077"

$bmPos = $cellRng.Start + 9
$bmRng = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRng)

# ---------------------------------------------------------------------------
# 2. "... cac truong day nghe nau an xuat hien ..." -> "... cac trung tam day
#    nau an xuat hien ..." and split the paragraph in two after "thoi gian? ".
# ---------------------------------------------------------------------------
$rngA = $d.Content
$rngA.Find.Execute("các trường dạy nghề nấu ăn xuất hiện", $false, $false, $false, $false, $false, $true, 1, $false, "các trung tâm dạy nấu ăn xuất hiện", 2) | Out-Null

$rngB = $d.Content
$rngB.Find.Execute("thời gian? ") | Out-Null
$splitPos = $rngB.End
$insAnchor = $d.Range($splitPos, $splitPos)
$insAnchor.InsertParagraphAfter()

# re-split "các " / "trung tâm" / " dạy nấu ..." into separate runs so the
# new phrase sits in its own run, matching the authored edit.
$rngC = $d.Content
$rngC.Find.Execute("trung tâm dạy nấu ăn xuất hiện") | Out-Null
$tt = $d.Range($rngC.Start, $rngC.Start + 9)
$tt.Font.Name = "Tahoma"

# ---------------------------------------------------------------------------
# 3. Merge "chu de. " + "Goi y cac khoa hoc noi bat nha" into a single run.
# ---------------------------------------------------------------------------
$rngD = $d.Content
$rngD.Find.Execute("chủ đề. Gợi ý các khóa học nổi bật nhấ") | Out-Null
$mergeRng = $d.Range($rngD.Start, $rngD.End)
$mergeRng.Text = $mergeRng.Text

Write-Output "done"
